$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changes from 46066 to 46070 for every data row (2..14)
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46070
}

# Rows 7..11 are refreshed with new/reordered records (columns A, B, G)
$updates = @(
    @{ Row = 7;  A = "A 46779-2025"; B = 45926;               G = 1.5 },
    @{ Row = 8;  A = "A 56948-2025"; B = 45978.64356481482;   G = 4.7 },
    @{ Row = 9;  A = "A 56917-2025"; B = 45978.58453703704;   G = 0.7 },
    @{ Row = 10; A = "A 31120-2023"; B = 45113;               G = 0.2 },
    @{ Row = 11; A = "A 64431-2023"; B = 45280;               G = 0.5 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.A
    $ws.Cells.Item($r, 2).Value2 = $u.B
    $ws.Cells.Item($r, 7).Value2 = $u.G
}
